$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 7713
$ws.Range("L2").Value = 6365
$ws.Range("L3").Value = 6868
$ws.Range("L4").Value = 1705
$ws.Range("K5").Value = 592
$ws.Range("L6").Value = 5650
$ws.Range("J7").Value = 29361
$ws.Range("K7").Value = 27586
$ws.Range("L7").Value = 20992

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 423
$ws.Range("L3").Value = 491
$ws.Range("L6").Value = 335
$ws.Range("L7").Value = 1387

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 147
$ws.Range("L6").Value = 104
$ws.Range("L7").Value = 459

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 259
$ws.Range("L7").Value = 944

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L4").Value = 46
$ws.Range("L6").Value = 208
$ws.Range("L7").Value = 806

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L6").Value = 103
$ws.Range("L7").Value = 412

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 111
$ws.Range("L7").Value = 363

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 188
$ws.Range("L5").Value = 74
$ws.Range("L7").Value = 667
$ws.Range("L8").Value = 1387
$ws.Range("L11").Value = 347
$ws.Range("L19").Value = 569
$ws.Range("L20").Value = 530
$ws.Range("L24").Value = 63
$ws.Range("L29").Value = 1174
$ws.Range("L30").Value = 93
$ws.Range("L33").Value = 944
$ws.Range("L37").Value = 806
$ws.Range("L44").Value = 144
$ws.Range("L46").Value = 50
$ws.Range("L47").Value = 148
$ws.Range("L51").Value = 259
$ws.Range("L54").Value = 454
$ws.Range("L55").Value = 224
$ws.Range("J63").Value = 237
$ws.Range("K63").Value = 181
$ws.Range("L63").Value = 63
$ws.Range("L64").Value = 133
$ws.Range("L65").Value = 412
$ws.Range("L70").Value = 38
$ws.Range("L73").Value = 165
$ws.Range("L79").Value = 582
$ws.Range("L83").Value = 459
$ws.Range("L84").Value = 201
$ws.Range("L85").Value = 1045
$ws.Range("L88").Value = 222
$ws.Range("L90").Value = 221
$ws.Range("L94").Value = 255
$ws.Range("L99").Value = 363
$ws.Range("L100").Value = 41
$ws.Range("J101").Value = 29361
$ws.Range("K101").Value = 27586
$ws.Range("L101").Value = 20992

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 69
$ws.Range("L7").Value = 201

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L2").Value = 84
$ws.Range("L3").Value = 113
$ws.Range("L7").Value = 454

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 355
$ws.Range("L6").Value = 285
$ws.Range("L7").Value = 1174

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L3").Value = 173
$ws.Range("L7").Value = 569

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 144

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L6").Value = 62
$ws.Range("L7").Value = 224

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("L2").Value = 17
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 180
$ws.Range("L3").Value = 187
$ws.Range("L6").Value = 156
$ws.Range("L7").Value = 582

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L3").Value = 41
$ws.Range("L7").Value = 133

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 168
$ws.Range("L7").Value = 530

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 212
$ws.Range("L7").Value = 667

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L2").Value = 61
$ws.Range("L7").Value = 255

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L3").Value = 51
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 148

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 130
$ws.Range("L7").Value = 347

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L2").Value = 59
$ws.Range("L7").Value = 165

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L2").Value = 62
$ws.Range("L7").Value = 188

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L3").Value = 75
$ws.Range("L7").Value = 222

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L4").Value = 19
$ws.Range("L7").Value = 221

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L6").Value = 56
$ws.Range("L7").Value = 259

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 314
$ws.Range("L3").Value = 435
$ws.Range("L7").Value = 1045
